$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.073.03"
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.338.04"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.39"
$ws.Range("E5").Value = "  -1.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.62"
$ws.Range("E6").Value = "  -2.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.40"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.960"
$ws.Range("E10").Value = "  +1.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.338.64"
$ws.Range("E11").Value = "  -1.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.99"
$ws.Range("E12").Value = "  +4.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.195"
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.18"
$ws.Range("E14").Value = "  +2.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.972.73"
$ws.Range("E15").Value = "  -1.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.968.25"
$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.13"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.340.20"
$ws.Range("E19").Value = "  -1.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.47"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.45"
$ws.Range("E22").Value = "  +9.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "494.74"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.442"
$ws.Range("E24").Value = "  -8.03%  "

$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.41"
$ws.Range("E25").Value = "  +1.36%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000183"
$ws.Range("E26").Value = "  -1.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "89.99"
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.93"
$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.515.93"
$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("E30").Value = "  -0.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.21"
$ws.Range("E31").Value = "  -2.53%  "

$ws.Range("E32").Value = "  +3.81%  "

$ws.Range("E33").Value = "  -2.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.32%  "

$ws.Range("E35").Value = "  -2.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.39"
$ws.Range("E36").Value = "  -4.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.530"
$ws.Range("E37").Value = "  -3.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "565.57"
$ws.Range("E38").Value = "  +5.53%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.39"
$ws.Range("E39").Value = "  -0.80%  "

$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.149"
$ws.Range("E41").Value = "  +0.55%  "

$ws.Range("E42").Value = "  -3.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.871"
$ws.Range("E43").Value = "  -4.17%  "

$ws.Range("E44").Value = "  -1.29%  "

$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.69"
$ws.Range("E45").Value = "  +1.97%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0415"
$ws.Range("E46").Value = "  +2.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.61"
$ws.Range("E47").Value = "  +6.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.43"
$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("E49").Value = "  +0.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.99"
$ws.Range("E50").Value = "  +1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.66"
$ws.Range("E51").Value = "  -1.84%  "
